$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 15 (pushes old row 15+ down by 2)
$ws.Range("A15:B16").EntireRow.Insert()

# Row 15: r4.4 "Prepare for Battle" overview entry
$ws.Range("A15").Value2 = "r4.4"
$ws.Range("B15").Value2 = "<Bold>r4.4 Prepare for Battle</Bold>`n<LineBreak/><LineBreak/>`n<InlineUIContainer><Button Content='r4.41' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Deployment <LineBreak/>`n<InlineUIContainer><Button Content='r4.42' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Hatches <LineBreak/>`n<InlineUIContainer><Button Content='r4.43' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Gun Load<LineBreak/>`n<InlineUIContainer><Button Content='r4.43' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Tank &amp; Turrent Orientation<LineBreak/>`n<InlineUIContainer><Button Content='r4.44' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Mark Loader Spot<LineBreak/>`n<InlineUIContainer><Button Content='r4.45' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Place <LineBreak/>`n<InlineUIContainer><Button Content='r4.46' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Place Use Control Markers<LineBreak/>`n<InlineUIContainer><Button Content='r4.47' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Mark Current Weather"
$ws.Rows.Item(15).RowHeight = 142.65

# Row 16: r4.41 "Deployment" entry
$ws.Range("A16").Value2 = "r4.41"
$ws.Range("B16").Value2 = "<Bold>r4.41 Deployment</Bold>`n<LineBreak/><LineBreak/>`nDetermine your tank&apos;s deployment from the `n<InlineUIContainer><Button Content='Deployment' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  Table.`n Your tank is marked as moving or hull down by having a counter placed on the Battle Board. `n<LineBreak/><LineBreak/>`nIf you tank is the lead tank, it is recorded on Notes section of the After Action Report (AAR) `n<InlineUIContainer><Button Content='r2.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>."
$ws.Rows.Item(16).RowHeight = 114.15

# Update selection to match the post-edit cursor position
$ws.Range("B12").Select()
